$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test account names/emails to the new "-test1" suite (row -> base name)
$names = @{
    2 = "ahostess-test1"
    3 = "bcohost-test1"
    4 = "guest1-test1"
    5 = "guest2-test1"
    6 = "guest3-test1"
}

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
}
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = "$($names[$r])@test.com"
}

# Update the active selection as recorded when the workbook was last saved
$ws.Range("E15").Select()
